$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: I2:L2
$ws.Range("I2").Value = -0.4695295520068083
$ws.Range("J2").Value = 0.2348751993419303
$ws.Range("K2").Value = -0.04957718332383978
$ws.Range("L2").Value = 2.767138876111711

# Row 20: I20:L20
$ws.Range("I20").Value = -0.6410027637727349
$ws.Range("J20").Value = 0.3214833734407471
$ws.Range("K20").Value = 0.02644294527821311
$ws.Range("L20").Value = 2.282637360416062
